$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures per commit "Updated cryptos list on Thu May 25 14:31:20 UTC 2023 with GitHub Actions"

$ws.Range("D2").Value = "26.450.18"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.805.46"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("E5").Value = "  -0.53%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "306.48"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.68%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4526"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.36%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3594"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.02%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "46.48"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07075"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.8886"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07808"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "19.44"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "1.853.47"
$ws.Range("E14").Value = "  +1.99%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.293"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.322"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "85.32"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").Value = "  -0.46%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008486"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "26.465.49"
$ws.Range("E21").Value = "  -0.44%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "14.21"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.972"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "2.068.89"
$ws.Range("E24").Value = "  +0.35%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "10.52"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.66%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.961"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "151.11"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "17.82"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.70%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.046"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.54%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "112.08"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.861"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.08692"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  +2.41%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.843"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +14.74%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.444"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7200"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.66%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.103"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  -0.36%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01935"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.05099"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +1.07%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.5120"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.42%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "6.783"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -3.59%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "8.009"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4663"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.66%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.51%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.988"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "100.48"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.573"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
